# Update "各级各类民办教育在校学生数" (enrollment) table: shift every year row
# up by one (2009 row drops off, each row now shows the following years
# figures) and append a new 2021 row at the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2010年
$ws.Range("A2").Value = '2010年'
$ws.Range("B2").Value = "'"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 306.9943
$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = 442.1129
$ws.Range("F2").Value = "'"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 1399.4694
$ws.Range("H2").Value = 442.1129
$ws.Range("I2").Value = 537.6255
$ws.Range("J2").Value = 230.0706
$ws.Range("K2").Value = 21.3403
$ws.Range("L2").Value = 260.3177
$ws.Range("M2").Value = 238.9774
$ws.Range("N2").Value = "'"
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").Value = 537.0649
$ws.Range("P2").Value = 195.6961
$ws.Range("Q2").Value = 476.6845
$ws.Range("R2").Value = 280.9884
$ws.Range("S2").Value = "'"
$ws.Range("S2").Style = "Normal"

# Row 3: 2011年
$ws.Range("A3").Value = '2011年'
$ws.Range("B3").Value = "'"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 269.2512
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 442.5616
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = 1694.209
$ws.Range("H3").Value = 442.5616
$ws.Range("I3").Value = 567.8255
$ws.Range("J3").Value = 234.9833
$ws.Range("K3").Value = 17.5738
$ws.Range("L3").Value = 267.4448
$ws.Range("M3").Value = 249.871
$ws.Range("N3").Value = "'"
$ws.Range("N3").Style = "Normal"
$ws.Range("O3").Value = 504.2345
$ws.Range("P3").Value = 193.2451
$ws.Range("Q3").Value = 505.0687
$ws.Range("R3").Value = 311.8236
$ws.Range("S3").Value = "'"
$ws.Range("S3").Style = "Normal"

# Row 4: 2012年
$ws.Range("A4").Value = '2012年'
$ws.Range("B4").Value = "'"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 240.88
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 451.4091
$ws.Range("F4").Value = "'"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = 1852.7444
$ws.Range("H4").Value = 451.4091
$ws.Range("I4").Value = 597.8535000000001
$ws.Range("J4").Value = 234.9575
$ws.Range("K4").Value = 16.249
$ws.Range("L4").Value = 278.3983
$ws.Range("M4").Value = 262.1493
$ws.Range("N4").Value = "'"
$ws.Range("N4").Style = "Normal"
$ws.Range("O4").Value = 475.8375
$ws.Range("P4").Value = 191.9358
$ws.Range("Q4").Value = 533.177
$ws.Range("R4").Value = 341.2257
$ws.Range("S4").Value = "'"
$ws.Range("S4").Style = "Normal"

# Row 5: 2013年
$ws.Range("A5").Value = '2013年'
$ws.Range("B5").Value = "'"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 207.9408
$ws.Range("D5").Value = "'"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = 462.3482
$ws.Range("F5").Value = "'"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = 1990.2536
$ws.Range("H5").Value = 462.3482
$ws.Range("I5").Value = 628.6015
$ws.Range("J5").Value = 231.6445
$ws.Range("K5").Value = 15.0206
$ws.Range("L5").Value = 275.8465
$ws.Range("M5").Value = 260.8259
$ws.Range("N5").Value = "'"
$ws.Range("N5").Style = "Normal"
$ws.Range("O5").Value = 439.5853
$ws.Range("P5").Value = 195.852
$ws.Range("Q5").Value = 557.5218
$ws.Range("R5").Value = 361.6363
$ws.Range("S5").Value = "'"
$ws.Range("S5").Style = "Normal"

# Row 6: 2014年
$ws.Range("A6").Value = '2014年'
$ws.Range("B6").Value = "'"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 189.5744
$ws.Range("D6").Value = "'"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = 487.0018
$ws.Range("F6").Value = "'"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = 2125.3781
$ws.Range("H6").Value = 486.995
$ws.Range("I6").Value = 674.1425
$ws.Range("J6").Value = 238.6542
$ws.Range("K6").Value = 13.6229
$ws.Range("L6").Value = 269.0625
$ws.Range("M6").Value = 255.4396
$ws.Range("N6").Value = "'"
$ws.Range("N6").Style = "Normal"
$ws.Range("O6").Value = 428.2286
$ws.Range("P6").Value = 212.2803
$ws.Range("Q6").Value = 587.1547
$ws.Range("R6").Value = 374.8336
$ws.Range("S6").Value = "'"
$ws.Range("S6").Style = "Normal"

# Row 7: 2015年
$ws.Range("A7").Value = '2015年'
$ws.Range("B7").Value = "'"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = 183.3658
$ws.Range("D7").Value = "'"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = 502.9261
$ws.Range("F7").Value = "'"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = 2302.4429
$ws.Range("H7").Value = 502.9229
$ws.Range("I7").Value = 713.8225
$ws.Range("J7").Value = 256.9644
$ws.Range("K7").Value = 11.7409
$ws.Range("L7").Value = 259.4243
$ws.Range("M7").Value = 247.6834
$ws.Range("N7").Value = "'"
$ws.Range("N7").Style = "Normal"
$ws.Range("O7").Value = 440.3302
$ws.Range("P7").Value = 227.5188
$ws.Range("Q7").Value = 610.9013
$ws.Range("R7").Value = 383.3316
$ws.Range("S7").Value = "'"
$ws.Range("S7").Style = "Normal"

# Row 8: 2016年
$ws.Range("A8").Value = '2016年'
$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 75.5553
$ws.Range("E8").Value = "'"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = 846.7954999999999
$ws.Range("G8").Value = 2437.6589
$ws.Range("H8").Value = "'"
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = 756.3291
$ws.Range("J8").Value = 279.0794
$ws.Range("K8").Value = "'"
$ws.Range("K8").Style = "Normal"
$ws.Range("L8").Value = "'"
$ws.Range("L8").Style = "Normal"
$ws.Range("M8").Value = "'"
$ws.Range("M8").Style = "Normal"
$ws.Range("N8").Value = 0.0051
$ws.Range("O8").Value = "'"
$ws.Range("O8").Style = "Normal"
$ws.Range("P8").Value = "'"
$ws.Range("P8").Style = "Normal"
$ws.Range("Q8").Value = "'"
$ws.Range("Q8").Style = "Normal"
$ws.Range("R8").Value = "'"
$ws.Range("R8").Style = "Normal"
$ws.Range("S8").Value = "'"
$ws.Range("S8").Style = "Normal"

# Row 9: 2017年
$ws.Range("A9").Value = '2017年'
$ws.Range("B9").Value = "'"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 222.2091
$ws.Range("D9").Value = "'"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = 577.6835
$ws.Range("F9").Value = 901.4369
$ws.Range("G9").Value = 2572.3434
$ws.Range("H9").Value = 577.6778
$ws.Range("I9").Value = 814.172
$ws.Range("J9").Value = 306.2608
$ws.Range("K9").Value = 10.1278
$ws.Range("L9").Value = 248.4741
$ws.Range("M9").Value = 238.3463
$ws.Range("N9").Value = "'"
$ws.Range("N9").Style = "Normal"
$ws.Range("O9").Value = "'"
$ws.Range("O9").Style = "Normal"
$ws.Range("P9").Value = 246.4354
$ws.Range("Q9").Value = 899.8533
$ws.Range("R9").Value = 404.7501
$ws.Range("S9").Value = "'"
$ws.Range("S9").Style = "Normal"

# Row 10: 2018年
$ws.Range("A10").Value = '2018年'
$ws.Range("B10").Value = "'"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = 8.5017
$ws.Range("E10").Value = "'"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = 910.0988
$ws.Range("G10").Value = 2639.7847
$ws.Range("H10").Value = "'"
$ws.Range("H10").Style = "Normal"
$ws.Range("I10").Value = 884.5746
$ws.Range("J10").Value = 328.2687
$ws.Range("K10").Value = "'"
$ws.Range("K10").Style = "Normal"
$ws.Range("L10").Value = "'"
$ws.Range("L10").Style = "Normal"
$ws.Range("M10").Value = "'"
$ws.Range("M10").Style = "Normal"
$ws.Range("N10").Value = 0.0101
$ws.Range("O10").Value = "'"
$ws.Range("O10").Style = "Normal"
$ws.Range("P10").Value = "'"
$ws.Range("P10").Style = "Normal"
$ws.Range("Q10").Value = "'"
$ws.Range("Q10").Style = "Normal"
$ws.Range("R10").Value = "'"
$ws.Range("R10").Style = "Normal"
$ws.Range("S10").Value = "'"
$ws.Range("S10").Style = "Normal"

# Row 11: 2019年
$ws.Range("A11").Value = '2019年'
$ws.Range("B11").Value = "'"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = 9.010300000000001
$ws.Range("E11").Value = "'"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = 1075.4715
$ws.Range("G11").Value = 2649.4401
$ws.Range("H11").Value = "'"
$ws.Range("H11").Style = "Normal"
$ws.Range("I11").Value = 944.9050999999999
$ws.Range("J11").Value = 359.6765
$ws.Range("K11").Value = "'"
$ws.Range("K11").Style = "Normal"
$ws.Range("L11").Value = "'"
$ws.Range("L11").Style = "Normal"
$ws.Range("M11").Value = "'"
$ws.Range("M11").Style = "Normal"
$ws.Range("N11").Value = 0.0091
$ws.Range("O11").Value = "'"
$ws.Range("O11").Style = "Normal"
$ws.Range("P11").Value = "'"
$ws.Range("P11").Style = "Normal"
$ws.Range("Q11").Value = "'"
$ws.Range("Q11").Style = "Normal"
$ws.Range("R11").Value = "'"
$ws.Range("R11").Style = "Normal"
$ws.Range("S11").Value = "'"
$ws.Range("S11").Style = "Normal"

# Row 12: 2020年
$ws.Range("A12").Value = '2020年'
$ws.Range("B12").Value = "'"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = 7.5528
$ws.Range("E12").Value = "'"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = 989.6564
$ws.Range("G12").Value = 2378.5506
$ws.Range("H12").Value = "'"
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").Value = 966.0348
$ws.Range("J12").Value = 401.2899
$ws.Range("K12").Value = "'"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = "'"
$ws.Range("L12").Style = "Normal"
$ws.Range("M12").Value = "'"
$ws.Range("M12").Style = "Normal"
$ws.Range("N12").Value = 0.0062
$ws.Range("O12").Value = "'"
$ws.Range("O12").Style = "Normal"
$ws.Range("P12").Value = "'"
$ws.Range("P12").Style = "Normal"
$ws.Range("Q12").Value = "'"
$ws.Range("Q12").Style = "Normal"
$ws.Range("R12").Value = "'"
$ws.Range("R12").Style = "Normal"
$ws.Range("S12").Value = "'"
$ws.Range("S12").Style = "Normal"

# Row 13: 2021年
$ws.Range("A13").Value = '2021年'
$ws.Range("B13").Value = "'"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 267.6329
$ws.Range("D13").Value = "'"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = 721.2729
$ws.Range("F13").Value = "'"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = 2312.0312
$ws.Range("H13").Value = "'"
$ws.Range("H13").Style = "Normal"
$ws.Range("I13").Value = 952.8280999999999
$ws.Range("J13").Value = 450.343
$ws.Range("K13").Value = "'"
$ws.Range("K13").Style = "Normal"
$ws.Range("L13").Value = "'"
$ws.Range("L13").Style = "Normal"
$ws.Range("M13").Value = "'"
$ws.Range("M13").Style = "Normal"
$ws.Range("N13").Value = "'"
$ws.Range("N13").Style = "Normal"
$ws.Range("O13").Value = 717.9759
$ws.Range("P13").Value = "'"
$ws.Range("P13").Style = "Normal"
$ws.Range("Q13").Value = "'"
$ws.Range("Q13").Style = "Normal"
$ws.Range("R13").Value = "'"
$ws.Range("R13").Style = "Normal"
$ws.Range("S13").Value = "'"
$ws.Range("S13").Style = "Normal"

